# Updates the cryptocurrency price list (Price and Volume(1h) columns)
# to the latest scraped values, and fixes the row data for two pairs of
# coins (InjectiveProtocol/TheGraph and ApeXProtocol/Stellar) whose rows
# were swapped in the source ranking so that each coin's name/link/price
# stays aligned with its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.809.08'
$ws.Range('E2').Value = '  +4.10%  '
$ws.Range('D3').Value = '3.440.81'
$ws.Range('E3').Value = '  +3.51%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Formula = "'579.22"
$ws.Range('E5').Value = '  +4.58%  '
$ws.Range('D6').Formula = "'184.19"
$ws.Range('E6').Value = '  +5.98%  '
$ws.Range('D7').Formula = "'0.631"
$ws.Range('E7').Value = '  +2.29%  '
$ws.Range('D8').Value = '3.436.66'
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Formula = "'0.172"
$ws.Range('E10').Value = '  +1.07%  '
$ws.Range('D11').Formula = "'0.645"
$ws.Range('E11').Value = '  +2.48%  '
$ws.Range('D12').Formula = "'56.10"
$ws.Range('E12').Value = '  +4.67%  '
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').Formula = "'9.40"
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').Value = '3.990.86'
$ws.Range('E15').Value = '  +3.52%  '
$ws.Range('D16').Formula = "'18.67"
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('D17').Value = '3.440.14'
$ws.Range('E17').Value = '  +3.50%  '
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '66.711.43'
$ws.Range('E19').Value = '  +2.32%  '
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').Formula = "'1.01"
$ws.Range('E21').Value = '  +3.09%  '
$ws.Range('D22').Formula = "'482.81"
$ws.Range('E22').Value = '  +6.58%  '
$ws.Range('D23').Formula = "'16.94"
$ws.Range('E23').Value = '  +23.31%  '
$ws.Range('D24').Formula = "'5.09"
$ws.Range('E24').Value = '  +2.97%  '
$ws.Range('D25').Formula = "'4.38"
$ws.Range('E25').Value = '  +7.49%  '
$ws.Range('D26').Formula = "'89.65"
$ws.Range('E26').Value = '  +3.35%  '
$ws.Range('E27').Value = '  +3.67%  '
$ws.Range('E28').Value = '  +2.93%  '
$ws.Range('D29').Formula = "'9.15"
$ws.Range('E29').Value = '  +6.70%  '
$ws.Range('D30').Formula = "'31.32"
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('D31').Formula = "'7.12"
$ws.Range('E31').Value = '  +8.72%  '
$ws.Range('D32').Formula = "'64.48"
$ws.Range('E32').Value = '  +6.51%  '
$ws.Range('E33').Value = '  +2.51%  '
$ws.Range('D34').Formula = "'592.32"
$ws.Range('E34').Value = '  +4.62%  '
$ws.Range('D35').Formula = "'0.111"
$ws.Range('E35').Value = '  +4.35%  '
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('E37').Value = '  +5.24%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Formula = "'36.38"
$ws.Range('E39').Value = '  +3.42%  '
$ws.Range('B40').Value = 'TheGraph'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D40').Formula = "'0.385"
$ws.Range('E40').Value = '  +5.33%  '
$ws.Range('D41').Value = '0.0₃0767'
$ws.Range('E41').Value = '  +4.67%  '
$ws.Range('D42').Value = '3.189.19'
$ws.Range('E42').Value = '  +4.16%  '
$ws.Range('D43').Formula = "'2.92"
$ws.Range('E43').Value = '  +5.53%  '
$ws.Range('E44').Value = '  +4.04%  '
$ws.Range('E45').Value = '  +4.97%  '
$ws.Range('E46').Value = '  +22.91%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Formula = "'0.135"
$ws.Range('E47').Value = '  +1.25%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Formula = "'3.20"
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Formula = "'8.70"
$ws.Range('E50').Value = '  +6.86%  '
$ws.Range('D51').Formula = "'139.73"
$ws.Range('E51').Value = '  -1.32%  '
